$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Value = "BonusPower"
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 960
$ws.Range("D30").Value = 96
$ws.Range("E30").Value = 100
$ws.Range("F30").Value = 200
$ws.Range("G30").Value = "win"
